$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the underlying raw counts for row 6 (2024-... data row)
$ws.Range("D6").Value = 75
$ws.Range("F6").Value = 484
$ws.Range("G6").Value = 37

$excel.CalculateFullRebuild()
